$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-09 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-10 Thursday", 2) | Out-Null
$d.Content.Find.Execute("12×75=900", $true, $false, $false, $false, $false, $true, 1, $false, "65×72=4680", 2) | Out-Null
$d.Content.Find.Execute("58×86=4988", $true, $false, $false, $false, $false, $true, 1, $false, "60×31=1860", 2) | Out-Null
$d.Content.Find.Execute("13×55=715", $true, $false, $false, $false, $false, $true, 1, $false, "32×87=2784", 2) | Out-Null
$d.Content.Find.Execute("95×57=5415", $true, $false, $false, $false, $false, $true, 1, $false, "18×72=1296", 2) | Out-Null
$d.Content.Find.Execute("93×31=2883", $true, $false, $false, $false, $false, $true, 1, $false, "13×61=793", 2) | Out-Null
$d.Content.Find.Execute("54×50=2700", $true, $false, $false, $false, $false, $true, 1, $false, "69×12=828", 2) | Out-Null
$d.Content.Find.Execute("48×23=1104", $true, $false, $false, $false, $false, $true, 1, $false, "27×41=1107", 2) | Out-Null
$d.Content.Find.Execute("19×86=1634", $true, $false, $false, $false, $false, $true, 1, $false, "37×51=1887", 2) | Out-Null
$d.Content.Find.Execute("31×11=341", $true, $false, $false, $false, $false, $true, 1, $false, "55×68=3740", 2) | Out-Null
$d.Content.Find.Execute("79×34=2686", $true, $false, $false, $false, $false, $true, 1, $false, "34×62=2108", 2) | Out-Null
$d.Content.Find.Execute("17×24=408", $true, $false, $false, $false, $false, $true, 1, $false, "59×47=2773", 2) | Out-Null
$d.Content.Find.Execute("55×69=3795", $true, $false, $false, $false, $false, $true, 1, $false, "38×72=2736", 2) | Out-Null
$d.Content.Find.Execute("29×42=1218", $true, $false, $false, $false, $false, $true, 1, $false, "37×91=3367", 2) | Out-Null
$d.Content.Find.Execute("47×56=2632", $true, $false, $false, $false, $false, $true, 1, $false, "30×30=900", 2) | Out-Null
$d.Content.Find.Execute("78×69=5382", $true, $false, $false, $false, $false, $true, 1, $false, "74×71=5254", 2) | Out-Null
$d.Content.Find.Execute("32×15=480", $true, $false, $false, $false, $false, $true, 1, $false, "90×45=4050", 2) | Out-Null
$d.Content.Find.Execute("18×88=1584", $true, $false, $false, $false, $false, $true, 1, $false, "59×47=2773", 2) | Out-Null
$d.Content.Find.Execute("91×58=5278", $true, $false, $false, $false, $false, $true, 1, $false, "54×44=2376", 2) | Out-Null
$d.Content.Find.Execute("64×65=4160", $true, $false, $false, $false, $false, $true, 1, $false, "23×66=1518", 2) | Out-Null
$d.Content.Find.Execute("64×45=2880", $true, $false, $false, $false, $false, $true, 1, $false, "53×91=4823", 2) | Out-Null
$d.Content.Find.Execute("17×11=187", $true, $false, $false, $false, $false, $true, 1, $false, "43×36=1548", 2) | Out-Null
$d.Content.Find.Execute("93×33=3069", $true, $false, $false, $false, $false, $true, 1, $false, "36×64=2304", 2) | Out-Null
$d.Content.Find.Execute("44×93=4092", $true, $false, $false, $false, $false, $true, 1, $false, "22×54=1188", 2) | Out-Null
$d.Content.Find.Execute("65×31=2015", $true, $false, $false, $false, $false, $true, 1, $false, "23×38=874", 2) | Out-Null
$d.Content.Find.Execute("76×89=6764", $true, $false, $false, $false, $false, $true, 1, $false, "55×98=5390", 2) | Out-Null
